$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new data row from DGS's 2021/09/24 report (next row after the
# existing last row, 84).
# Column A holds the report date as text (matching the rest of the column),
# so force a text format before assigning it, then restore the date display
# format used by the rest of the column so the cell's style matches its
# neighbours.
$ws.Range("A85").NumberFormat = "@"
$ws.Range("A85").Value = "2021/09/24"
$ws.Range("A85").NumberFormat = "yyyy/mm/dd"

$ws.Range("B85").Value = 127.3
$ws.Range("C85").Value = 129.7
$ws.Range("D85").Value = 0.83
$ws.Range("E85").Value = 0.82

# Move the active selection to the next empty row, as happens in Excel/Calc
# after entering data in the previous row.
$ws.Range("A86").Select()
